$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.750.39'
$ws.Range('E2').Value = '  +1.03%  '
$ws.Range('D3').Value = '3.282.26'
$ws.Range('E3').Value = '  +5.12%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '599.13'
$ws.Range('E5').Value = '  +1.17%  '
$ws.Range('D6').Value = '143.31'
$ws.Range('E6').Value = '  +5.25%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.277.81'
$ws.Range('E8').Value = '  +5.20%  '
$ws.Range('D9').Value = '0.522'
$ws.Range('E9').Value = '  +0.58%  '
$ws.Range('E10').Value = '  +2.27%  '
$ws.Range('D11').Value = '5.44'
$ws.Range('E11').Value = '  +1.58%  '
$ws.Range('D12').Value = '0.472'
$ws.Range('E12').Value = '  +3.07%  '
$ws.Range('D13').Value = '0.0000248'
$ws.Range('E13').Value = '  -0.14%  '
$ws.Range('D14').Value = '34.89'
$ws.Range('E14').Value = '  +2.46%  '
$ws.Range('D15').Value = '3.818.34'
$ws.Range('E15').Value = '  +5.16%  '
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('D17').Value = '3.278.37'
$ws.Range('E17').Value = '  +4.94%  '
$ws.Range('D18').Value = '63.800.07'
$ws.Range('E18').Value = '  +1.08%  '
$ws.Range('D19').Value = '6.90'
$ws.Range('E19').Value = '  +3.18%  '
$ws.Range('D20').Value = '481.36'
$ws.Range('E20').Value = '  +1.26%  '
$ws.Range('D21').Value = '14.28'
$ws.Range('E21').Value = '  +0.90%  '
$ws.Range('D22').Value = '0.743'
$ws.Range('E22').Value = '  +6.67%  '
$ws.Range('D23').Value = '8.04'
$ws.Range('E23').Value = '  +5.21%  '
$ws.Range('D24').Value = '13.52'
$ws.Range('E24').Value = '  +4.04%  '
$ws.Range('D25').Value = '84.59'
$ws.Range('E25').Value = '  -2.73%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  +2.02%  '
$ws.Range('D28').Value = '7.30'
$ws.Range('E28').Value = '  +2.12%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '8.28'
$ws.Range('E29').Value = '  +3.53%  '
$ws.Range('B30').Value = 'FirstDigitalUSD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('D31').Value = '2.17'
$ws.Range('E31').Value = '  +5.79%  '
$ws.Range('D32').Value = '28.15'
$ws.Range('E32').Value = '  +3.90%  '
$ws.Range('E33').Value = '  +0.79%  '
$ws.Range('D34').Value = '2.57'
$ws.Range('E34').Value = '  +1.33%  '
$ws.Range('D35').Value = '1.10'
$ws.Range('E35').Value = '  +2.28%  '
$ws.Range('D36').Value = '6.00'
$ws.Range('E36').Value = '  +2.64%  '
$ws.Range('D37').Value = '53.12'
$ws.Range('E37').Value = '  +2.03%  '
$ws.Range('D38').Value = '0.0₃0737'
$ws.Range('E38').Value = '  +3.37%  '
$ws.Range('D39').Value = '0.0398'
$ws.Range('E39').Value = '  +2.12%  '
$ws.Range('D40').Value = '427.53'
$ws.Range('E40').Value = '  +1.41%  '
$ws.Range('D41').Value = '3.019.26'
$ws.Range('E41').Value = '  +5.45%  '
$ws.Range('D42').Value = '8.47'
$ws.Range('E42').Value = '  +2.61%  '
$ws.Range('D43').Value = '2.79'
$ws.Range('E43').Value = '  +3.32%  '
$ws.Range('E44').Value = '  -3.84%  '
$ws.Range('D45').Value = '0.270'
$ws.Range('E45').Value = '  +5.14%  '
$ws.Range('D46').Value = '2.25'
$ws.Range('E46').Value = '  +6.88%  '
$ws.Range('D47').Value = '26.28'
$ws.Range('E47').Value = '  +3.45%  '
$ws.Range('D49').Value = '2.35'
$ws.Range('E49').Value = '  +3.07%  '
$ws.Range('E50').Value = '  +1.43%  '
$ws.Range('D51').Value = '123.00'
$ws.Range('E51').Value = '  +3.60%  '
